# Auto-generated edit script: applies market-price data refresh to Sargatanas_Profits workbook
# (H=currentAveragePrice, I=currentAveragePriceNQ, J=currentAveragePriceHQ,
#  K=LevePriceNQ, L=LevePriceHQ, M=LeveProfitNQ, N=LeveProfitHQ)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 850
$ws.Range("I19").Value = 800
$ws.Range("K19").Value = 800
$ws.Range("M19").Value = -625
$ws.Range("H28").Value = 409
$ws.Range("I28").Value = 426.41666
$ws.Range("K28").Value = 426.41666
$ws.Range("M28").Value = 58.58334000000002
$ws.Range("H132").Value = 3932.1538
$ws.Range("I132").Value = 3932.1538
$ws.Range("K132").Value = 11796.4614
$ws.Range("M132").Value = -9266.4614
$ws.Range("H141").Value = 1704.5652
$ws.Range("I141").Value = 1507.1111
$ws.Range("K141").Value = 4521.3333
$ws.Range("M141").Value = 658.6666999999998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4200.353
$ws.Range("J61").Value = 6217.8335
$ws.Range("L61").Value = 6217.8335
$ws.Range("N61").Value = -6641.8335
$ws.Range("H136").Value = 4200.353
$ws.Range("J136").Value = 6217.8335
$ws.Range("L136").Value = 18653.5005
$ws.Range("N136").Value = -23753.5005

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H55").Value = 59376
$ws.Range("J55").Value = 59376
$ws.Range("L55").Value = 59376
$ws.Range("N55").Value = -59922
$ws.Range("H64").Value = 22223382
$ws.Range("I64").Value = 47619900
$ws.Range("J64").Value = 1428.5
$ws.Range("K64").Value = 47619900
$ws.Range("L64").Value = 1428.5
$ws.Range("M64").Value = -47619675
$ws.Range("N64").Value = -1878.5
$ws.Range("H67").Value = 22223382
$ws.Range("I67").Value = 47619900
$ws.Range("J67").Value = 1428.5
$ws.Range("K67").Value = 47619900
$ws.Range("L67").Value = 1428.5
$ws.Range("M67").Value = -47619120
$ws.Range("N67").Value = -2988.5
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()
$ws.Range("H105").Value = 2874.75
$ws.Range("I105").Value = 2499.182
$ws.Range("K105").Value = 2499.182
$ws.Range("M105").Value = -752.1819999999998
$ws.Range("H128").Value = 3840
$ws.Range("I128").Value = 3840
$ws.Range("K128").Value = 11520
$ws.Range("M128").Value = -9030
$ws.Range("H134").Value = 3938.1973
$ws.Range("I134").Value = 1976.7627
$ws.Range("J134").Value = 10745.529
$ws.Range("K134").Value = 5930.2881
$ws.Range("L134").Value = 32236.587
$ws.Range("M134").Value = -3395.2881
$ws.Range("N134").Value = -37306.587

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6222.6816
$ws.Range("I31").Value = 2716.6667
$ws.Range("J31").Value = 12358.208
$ws.Range("K31").Value = 2716.6667
$ws.Range("L31").Value = 12358.208
$ws.Range("M31").Value = -2421.6667
$ws.Range("N31").Value = -12948.208
$ws.Range("H34").Value = 6222.6816
$ws.Range("I34").Value = 2716.6667
$ws.Range("J34").Value = 12358.208
$ws.Range("K34").Value = 2716.6667
$ws.Range("L34").Value = 12358.208
$ws.Range("M34").Value = -2514.6667
$ws.Range("N34").Value = -12762.208
$ws.Range("H62").Value = 6717.909
$ws.Range("I62").Value = 6740.875
$ws.Range("J62").Value = 6656.6665
$ws.Range("K62").Value = 6740.875
$ws.Range("L62").Value = 6656.6665
$ws.Range("M62").Value = -6116.875
$ws.Range("N62").Value = -7904.6665
$ws.Range("H65").Value = 6717.909
$ws.Range("I65").Value = 6740.875
$ws.Range("J65").Value = 6656.6665
$ws.Range("K65").Value = 33704.375
$ws.Range("L65").Value = 33283.3325
$ws.Range("M65").Value = -30584.375
$ws.Range("N65").Value = -39523.3325
$ws.Range("H134").Value = 9986.571
$ws.Range("I134").Value = 12646
$ws.Range("K134").Value = 37938
$ws.Range("M134").Value = -35403

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 16286002
$ws.Range("I4").Value = 19194208
$ws.Range("J4").Value = 49.8
$ws.Range("K4").Value = 57582624
$ws.Range("L4").Value = 149.4
$ws.Range("M4").Value = -57582512
$ws.Range("N4").Value = -373.4
$ws.Range("H12").Value = 2174610.8
$ws.Range("I12").Value = 934.75
$ws.Range("J12").Value = 3333904.5
$ws.Range("K12").Value = 2804.25
$ws.Range("L12").Value = 10001713.5
$ws.Range("M12").Value = -2631.25
$ws.Range("N12").Value = -10002059.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 2169.7222
$ws.Range("I97").Value = 1896.5
$ws.Range("J97").Value = 2511.25
$ws.Range("K97").Value = 1896.5
$ws.Range("L97").Value = 2511.25
$ws.Range("M97").Value = -1400.5
$ws.Range("N97").Value = -3503.25
$ws.Range("H126").Value = 2492.0476
$ws.Range("I126").Value = 2355.6667
$ws.Range("J126").Value = 2594.3333
$ws.Range("K126").Value = 7067.000100000001
$ws.Range("L126").Value = 7782.999899999999
$ws.Range("M126").Value = -4597.000100000001
$ws.Range("N126").Value = -12722.9999
$ws.Range("H132").Value = 2084.9333
$ws.Range("I132").Value = 2042.7
$ws.Range("J132").Value = 2422.8
$ws.Range("K132").Value = 6128.1
$ws.Range("L132").Value = 7268.400000000001
$ws.Range("M132").Value = -3598.1
$ws.Range("N132").Value = -12328.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2901
$ws.Range("I22").Value = 1333
$ws.Range("K22").Value = 1333
$ws.Range("M22").Value = -1038
$ws.Range("H27").Value = 2901
$ws.Range("I27").Value = 1333
$ws.Range("K27").Value = 1333
$ws.Range("M27").Value = -1226
$ws.Range("H93").Value = 9572.857
$ws.Range("I93").Value = 11000.75
$ws.Range("J93").Value = 7669
$ws.Range("K93").Value = 11000.75
$ws.Range("L93").Value = 7669
$ws.Range("M93").Value = -9752.75
$ws.Range("N93").Value = -10165
$ws.Range("H100").Value = 3543.7144
$ws.Range("I100").Value = 3199.75
$ws.Range("J100").Value = 4002.3333
$ws.Range("K100").Value = 3199.75
$ws.Range("L100").Value = 4002.3333
$ws.Range("M100").Value = -2658.75
$ws.Range("N100").Value = -5084.3333

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 10674.444
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").ClearContents()
$ws.Range("H62").Value = 12000
$ws.Range("I62").Value = 12000
$ws.Range("K62").Value = 12000
$ws.Range("M62").Value = -11376
$ws.Range("H65").Value = 12000
$ws.Range("I65").Value = 12000
$ws.Range("K65").Value = 60000
$ws.Range("M65").Value = -56880
$ws.Range("H81").Value = 20018888
$ws.Range("I81").Value = 1648
$ws.Range("J81").Value = 50044748
$ws.Range("K81").Value = 3296
$ws.Range("L81").Value = 100089496
$ws.Range("M81").Value = -2235
$ws.Range("N81").Value = -100091618
$ws.Range("H84").Value = 20018888
$ws.Range("I84").Value = 1648
$ws.Range("J84").Value = 50044748
$ws.Range("K84").Value = 16480
$ws.Range("L84").Value = 500447480
$ws.Range("M84").Value = -11176
$ws.Range("N84").Value = -500458088
$ws.Range("H107").Value = 1485
$ws.Range("I107").Value = 1601.5
$ws.Range("K107").Value = 4804.5
$ws.Range("M107").Value = -2884.5
$ws.Range("H113").Value = 2941.5
$ws.Range("I113").Value = 2550
$ws.Range("K113").Value = 7650
$ws.Range("M113").Value = -5480
$ws.Range("H136").Value = 47623544
$ws.Range("J136").Value = 7706.364
$ws.Range("L136").Value = 23119.092
$ws.Range("N136").Value = -28219.092
